# Add "Total Category Level Data" / cat_ rows for the new volume-pack
# variables to the taxonomy sheet (rows 2919-2964), matching the new
# sharedStrings + sheet2 rows from the commit "different versions of
# models, incl. TotalCategory variables".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("taxonomy")
$ws.Activate()

# row, D (classification), E (abbreviation2)
$rows = @(
  @(2919, "volume multiples glass pack", "vol_multiples_glass_single"),
  @(2920, "volume multiples glass pack", "vol_multiples_glass_2_pack"),
  @(2921, "volume multiples glass pack", "vol_multiples_glass_3_pack"),
  @(2922, "volume multiples glass pack", "vol_multiples_glass_4_pack"),
  @(2923, "volume multiples glass pack", "vol_multiples_glass_6_pack"),
  @(2924, "volume multiples glass pack", "vol_multiples_glass_8_pack"),
  @(2925, "volume multiples glass pack", "vol_multiples_glass_10_pack"),
  @(2926, "volume multiples glass pack", "vol_multiples_glass_12_pack"),
  @(2927, "volume multiples glass pack", "vol_multiples_glass_15_pack"),
  @(2928, "volume multiples glass pack", "vol_multiples_glass_18_pack"),
  @(2929, "volume multiples glass pack", "vol_multiples_glass_20_pack"),
  @(2930, "volume multiples glass pack", "vol_multiples_glass_24_pack"),
  @(2931, "volume multiples can pack", "vol_multiples_can_single"),
  @(2932, "volume multiples can pack", "vol_multiples_can_2_pack"),
  @(2933, "volume multiples can pack", "vol_multiples_can_4_pack"),
  @(2934, "volume multiples can pack", "vol_multiples_can_6_pack"),
  @(2935, "volume multiples can pack", "vol_multiples_can_8_pack"),
  @(2936, "volume multiples can pack", "vol_multiples_can_10_pack"),
  @(2937, "volume multiples can pack", "vol_multiples_can_12_pack"),
  @(2938, "volume multiples can pack", "vol_multiples_can_15_pack"),
  @(2939, "volume multiples can pack", "vol_multiples_can_18_pack"),
  @(2940, "volume multiples can pack", "vol_multiples_can_20_pack"),
  @(2941, "volume multiples can pack", "vol_multiples_can_24_pack"),
  @(2942, "volume impulse glass pack", "vol_impulse_glass_single"),
  @(2943, "volume impulse glass pack", "vol_impulse_glass_2_pack"),
  @(2944, "volume impulse glass pack", "vol_impulse_glass_3_pack"),
  @(2945, "volume impulse glass pack", "vol_impulse_glass_4_pack"),
  @(2946, "volume impulse glass pack", "vol_impulse_glass_5_pack"),
  @(2947, "volume impulse glass pack", "vol_impulse_glass_6_pack"),
  @(2948, "volume impulse glass pack", "vol_impulse_glass_10_pack"),
  @(2949, "volume impulse glass pack", "vol_impulse_glass_12_pack"),
  @(2950, "volume impulse glass pack", "vol_impulse_glass_15_pack"),
  @(2951, "volume impulse glass pack", "vol_impulse_glass_8_pack"),
  @(2952, "volume impulse glass pack", "vol_impulse_glass_18_pack"),
  @(2953, "volume impulse glass pack", "vol_impulse_glass_20_pack"),
  @(2954, "volume multiples can pack", "vol_impulse_can_24_pack"),
  @(2955, "volume multiples can pack", "vol_impulse_can_single"),
  @(2956, "volume multiples can pack", "vol_impulse_can_4_pack"),
  @(2957, "volume multiples can pack", "vol_impulse_can_6_pack"),
  @(2958, "volume multiples can pack", "vol_impulse_can_10_pack"),
  @(2959, "volume multiples can pack", "vol_impulse_can_12_pack"),
  @(2960, "volume multiples can pack", "vol_impulse_can_15_pack"),
  @(2961, "volume multiples can pack", "vol_impulse_can_8_pack"),
  @(2962, "volume multiples can pack", "vol_impulse_can_18_pack"),
  @(2963, "volume multiples can pack", "vol_impulse_can_20_pack"),
  @(2964, "volume multiples can pack", "vol_impulse_can_24_pack")
)

foreach ($row in $rows) {
  $r = $row[0]
  $d = $row[1]
  $e = $row[2]
  $ws.Cells.Item($r, 2).Value = "Total Category Level Data"
  $ws.Cells.Item($r, 3).Value = "cat_"
  $ws.Cells.Item($r, 4).Value = $d
  $ws.Cells.Item($r, 5).Value = $e
  $ws.Cells.Item($r, 8).Value = "market"
}

# Column A is the shared concat formula used throughout the sheet:
# variable_name = category_abbrev & abbreviation2 & (optional abbreviation3)
$ws.Range("A2919:A2964").Formula = "=C2919&E2919&G2919"

# Match the new selection left by the edit (sheet view stays on taxonomy,
# frozen header pane unchanged).
$ws.Range("A2946").Select()

Write-Output "Added taxonomy rows 2919:2964 (Total Category Level Data / cat_)"
